$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new timestamped reading for 2026/01/14 (水) was appended to the log,
# which gets inserted right after the existing 2026/01/14 entries (row 647),
# pushing every subsequent row down by one.
$ws.Rows(647).Copy()
$ws.Rows(648).Insert()

# Row 648 now duplicates row 647's "2026/01/14" / "水" text exactly as
# needed; only the time-of-day (C) and ranking (D) values differ.
$ws.Range("C648").Value = 14
$ws.Range("D648").Value = 201
